# Update "want to go" counts (column F) that changed between site crawls.
# Both the "展览" sheet and the "全部类型" sheet contain the same events
# (on different rows), so the same update is applied on each.

$wb = $excel.ActiveWorkbook

# Sheet "展览" (rows 3,4,5,7,9,10,11,12,13)
$ws1 = $wb.Worksheets.Item("展览")
$ws1.Range("F3").Value  = 1021
$ws1.Range("F4").Value  = 165
$ws1.Range("F5").Value  = 2765
$ws1.Range("F7").Value  = 217
$ws1.Range("F9").Value  = 119
$ws1.Range("F10").Value = 58
$ws1.Range("F11").Value = 61
$ws1.Range("F12").Value = 2579
$ws1.Range("F13").Value = 734

# Sheet "全部类型" (rows 4,5,6,8,11,12,13,14,15)
$ws4 = $wb.Worksheets.Item("全部类型")
$ws4.Range("F4").Value  = 1021
$ws4.Range("F5").Value  = 165
$ws4.Range("F6").Value  = 2765
$ws4.Range("F8").Value  = 217
$ws4.Range("F11").Value = 119
$ws4.Range("F12").Value = 58
$ws4.Range("F13").Value = 61
$ws4.Range("F14").Value = 2579
$ws4.Range("F15").Value = 734
